# Fruta / hortaliza, semanal
# A new weekly record is inserted as row 9, pushing all subsequent data rows
# down by one (old row 9 becomes row 10, ..., old row 50 becomes row 51).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 9, shifting existing rows 9-50 down to 10-51.
$ws.Rows.Item(9).Insert()

# Populate the new row 9 with this week's record for the same
# Vega Modelo de Temuco / Tuna series.
$ws.Range("A9").Value = 10
$ws.Range("B9").Value = "Vega Modelo de Temuco"
$ws.Range("C9").Value = "La Araucanía"
$ws.Range("D9").Value = 44630
$ws.Range("E9").Value = 9
$ws.Range("F9").Value = "Fruta"
$ws.Range("G9").Value = 100107
$ws.Range("H9").Value = "Otros"
$ws.Range("I9").Value = 100107011
$ws.Range("J9").Value = "Tuna"
$ws.Range("K9").Value = "Sin especificar"
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 100
$ws.Range("N9").Value = 17000
$ws.Range("O9").Value = 17000
$ws.Range("P9").Value = 17000
$ws.Range("Q9").Value = "$/caja 16 kilos"
$ws.Range("R9").Value = "Provincia de Los Andes"
$ws.Range("S9").Value = 1062
$ws.Range("T9").Value = 16
